$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "True"
